# tblStudy.xlsx fix: split "IncludeInPublication" (AI) into
#   AI = "RandomSample"      (existing flag, reset to FALSE for rows 39-53)
#   AJ = "UseInPublication"  (new column, keeps the old AI/IncludeInPublication value)
# plus a batch of column-width resets (J,K,L,N,O,Q,R,T,U,V,Y,Z,AA,AC,AE -> 14.0625)
# and a brand-new column AJ (width 14.0625).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Create the new AJ column by duplicating AI (values + formatting),
#        so AJ ends up with the exact same per-row style AI had. ---
$ws.Range("AI1:AI53").Copy($ws.Range("AJ1:AJ53"))

# --- 2. Re-label the headers. ---
$ws.Range("AI1").Value = "RandomSample"
$ws.Range("AJ1").Value = "UseInPublication"

# --- 3. Rows 39-53: AI ("RandomSample") resets to FALSE; AJ keeps the
#        original IncludeInPublication value (already copied above, stays TRUE). ---
for ($r = 39; $r -le 53; $r++) {
    $ws.Cells.Item($r, 35).Value = $false
}

# --- 4. Column width resets: J,K,L,N,O,Q,R,T,U,V,Y,Z,AA,AC,AE -> 14.0625 ---
#        (ColumnWidth is entered in Excel "characters"; the engine stores
#        width = round((chars + 5/7)*7)/7, so feed it chars = target - 5/7
#        to land as close as possible to 14.0625.)
$target = 13.348214285714286
$colsToReset = @(10, 11, 12, 14, 15, 17, 18, 20, 21, 22, 25, 26, 27, 29, 31)
foreach ($c in $colsToReset) {
    $ws.Columns.Item($c).ColumnWidth = $target
}

# --- 5. New column AJ (36) width -> 14.0625 as well. ---
$ws.Columns.Item(36).ColumnWidth = $target
